$d = $word.ActiveDocument
$dash = [char]0x2013

# ---------------------------------------------------------------------------
# 1) Title: "h-Mn part plan" -> "Strain induced coherent dynamics of Mn-doped
#    positively charged quantum dots"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "h-Mn part plan", $true, $false, $false, $false, $false, $true, 1, $false,
    "Strain induced coherent dynamics of Mn-doped positively charged quantum dots", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) "I.1 " + en-dash + " Energy structure" -> "I.2 " + en-dash + " Energy structure"
#    (the section was renumbered from I.1 to I.2)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "I.1 " + $dash + " Energy structure", $true, $false, $false, $false, $false, $true, 1, $false,
    "I.2 " + $dash + " Energy structure", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 3) "I.2 " + en-dash + " Optical " -> "I.3 " + en-dash + " Optical "
#    (the following section was renumbered from I.2 to I.3)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "I.2 " + $dash + " Optical ", $true, $false, $false, $false, $false, $true, 1, $false,
    "I.3 " + $dash + " Optical ", 2
) | Out-Null

# ---------------------------------------------------------------------------
# 4) Hyperlink display text "XplusMnR" + "e" + "s.pptx" collapses into a single
#    run "XplusMnRes.pptx" (no textual change, only a run-merge triggered by
#    touching the range).
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("XplusMnRes.pptx") | Out-Null
$hStart = $rng.Start
$hEnd = $rng.End
$tmp = $d.Range($hEnd, $hEnd)
$tmp.InsertAfter("ZZ")
$tmp2 = $d.Range($hEnd, $hEnd + 2)
$tmp2.Delete()

# ---------------------------------------------------------------------------
# 5) "|3, +1> and |3, -1>, and |2, +1> and |2, -1> coupled by E" collapses
#    from many small runs into a single run (no textual change).
# ---------------------------------------------------------------------------
$target5 = "|3, +1> and |3, -1>, and |2, +1> and |2, -1> coupled by E"
$rng = $d.Content
$rng.Find.Execute($target5) | Out-Null
$s5 = $rng.Start
$e5 = $rng.End
$tmp = $d.Range($e5, $e5)
$tmp.InsertAfter("ZZ")
$tmp2 = $d.Range($e5, $e5 + 2)
$tmp2.Delete()

# ---------------------------------------------------------------------------
# 6) "Experiment configuration |3, +1> " + " + " collapses into a single run
#    "Experiment configuration |3, +1>  + " (no textual change).
# ---------------------------------------------------------------------------
$target6 = "Experiment configuration |3, +1>  + "
$rng = $d.Content
$rng.Find.Execute($target6) | Out-Null
$s6 = $rng.Start
$e6 = $rng.End
$tmp = $d.Range($e6, $e6)
$tmp.InsertAfter("ZZ")
$tmp2 = $d.Range($e6, $e6 + 2)
$tmp2.Delete()

# ---------------------------------------------------------------------------
# 7) Final block: an extra empty paragraph is inserted after the "Schema of
#    the QD..." paragraph, the trailing "_GoBack" bookmark moves from the end
#    of that paragraph to between "Fig.15" and ": Polarization rate
#    evolution..." in the last paragraph, and "Fig.12" becomes "Fig.15"
#    there.
# ---------------------------------------------------------------------------
$lastTarget = "Fig.12: Polarization rate evolution in B(x and z) and simulation"
$d.Content.Find.Execute($lastTarget) | Out-Null

# Find the paragraph that currently holds this text and insert a new empty
# paragraph right before it (mirrors the pre-existing empty paragraph just
# above it).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq $lastTarget) {
        $para.Range.InsertParagraphBefore()
        break
    }
}

# Remove the "_GoBack" bookmark from its old location.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Locate the (now shifted) paragraph again and rewrite "Fig.12" -> "Fig.15".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq $lastTarget) {
        $figRange = $d.Range($para.Range.Start, $para.Range.Start + 6)
        $figRange.Text = "Fig.15"
        break
    }
}

# Re-insert the "_GoBack" bookmark right after the new "Fig.15" text.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -eq "Fig.15: Polarization rate evolution in B(x and z) and simulation") {
        $bmPos = $para.Range.Start + 6
        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
        break
    }
}

Write-Output "done"
